$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 318.5625
$ws.Range("I39").Value = 296.46667
$ws.Range("J39").Value = 650
$ws.Range("K39").Value = 889.4000100000001
$ws.Range("L39").Value = 1950
$ws.Range("M39").Value = -593.4000100000001
$ws.Range("N39").Value = -2542
$ws.Range("H43").Value = 2090.3333
$ws.Range("J43").Value = 2101.625
$ws.Range("L43").Value = 2101.625
$ws.Range("N43").Value = -2239.625
$ws.Range("H62").Value = 10550
$ws.Range("I62").Value = 10035.462
$ws.Range("K62").Value = 10035.462
$ws.Range("M62").Value = -9411.462
$ws.Range("H65").Value = 10550
$ws.Range("I65").Value = 10035.462
$ws.Range("K65").Value = 50177.31
$ws.Range("M65").Value = -47057.31
$ws.Range("H70").Value = 41686.516
$ws.Range("I70").Value = 10314.286
$ws.Range("J70").Value = 50836.75
$ws.Range("K70").Value = 30942.858
$ws.Range("L70").Value = 152510.25
$ws.Range("M70").Value = -30672.858
$ws.Range("N70").Value = -153050.25
$ws.Range("H73").Value = 41686.516
$ws.Range("I73").Value = 10314.286
$ws.Range("J73").Value = 50836.75
$ws.Range("K73").Value = 30942.858
$ws.Range("L73").Value = 152510.25
$ws.Range("M73").Value = -30006.858
$ws.Range("N73").Value = -154382.25
$ws.Range("H74").Value = 6959.875
$ws.Range("I74").Value = 6913.3335
$ws.Range("J74").Value = 7099.5
$ws.Range("K74").Value = 6913.3335
$ws.Range("L74").Value = 7099.5
$ws.Range("M74").Value = -5977.3335
$ws.Range("N74").Value = -8971.5
$ws.Range("H77").Value = 6959.875
$ws.Range("I77").Value = 6913.3335
$ws.Range("J77").Value = 7099.5
$ws.Range("K77").Value = 34566.6675
$ws.Range("L77").Value = 35497.5
$ws.Range("M77").Value = -29886.6675
$ws.Range("N77").Value = -44857.5
$ws.Range("H86").Value = 2547.2
$ws.Range("I86").Value = 2997.7144
$ws.Range("J86").Value = 1496
$ws.Range("K86").Value = 2997.7144
$ws.Range("L86").Value = 1496
$ws.Range("M86").Value = -1874.7144
$ws.Range("N86").Value = -3742
$ws.Range("H89").Value = 2547.2
$ws.Range("I89").Value = 2997.7144
$ws.Range("J89").Value = 1496
$ws.Range("K89").Value = 14988.572
$ws.Range("L89").Value = 7480
$ws.Range("M89").Value = -9372.572
$ws.Range("N89").Value = -18712
$ws.Range("H106").Value = 3398.5
$ws.Range("I106").Value = 3398.5
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3398.5
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -2767.5
$ws.Range("N106").ClearContents()
$ws.Range("H116").Value = 3971.875
$ws.Range("I116").Value = 3700
$ws.Range("J116").Value = 4425
$ws.Range("K116").Value = 3700
$ws.Range("L116").Value = 4425
$ws.Range("M116").Value = -258
$ws.Range("N116").Value = -11309
$ws.Range("H129").Value = 1138.2632
$ws.Range("I129").Value = 750.5
$ws.Range("K129").Value = 2251.5
$ws.Range("M129").Value = 2748.5
$ws.Range("H132").Value = 17905.047
$ws.Range("I132").Value = 17905.047
$ws.Range("K132").Value = 53715.141
$ws.Range("M132").Value = -51185.141
$ws.Range("H138").Value = 4414.4365
$ws.Range("I138").Value = 3709.625
$ws.Range("J138").Value = 4703.59
$ws.Range("K138").Value = 11128.875
$ws.Range("L138").Value = 14110.77
$ws.Range("M138").Value = -5988.875
$ws.Range("N138").Value = -24390.77

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6738.3555
$ws.Range("I32").Value = 6374.2617
$ws.Range("K32").Value = 6374.2617
$ws.Range("M32").Value = -6087.2617
$ws.Range("H61").Value = 2691.2808
$ws.Range("I61").Value = 2148.889
$ws.Range("J61").Value = 4725.25
$ws.Range("K61").Value = 2148.889
$ws.Range("L61").Value = 4725.25
$ws.Range("M61").Value = -1936.889
$ws.Range("N61").Value = -5149.25
$ws.Range("H74").Value = 184063.19
$ws.Range("I74").Value = 262242.6
$ws.Range("K74").Value = 262242.6
$ws.Range("M74").Value = -261368.6
$ws.Range("H77").Value = 184063.19
$ws.Range("I77").Value = 262242.6
$ws.Range("K77").Value = 1311213
$ws.Range("M77").Value = -1306845
$ws.Range("H109").Value = 82763.5
$ws.Range("J109").Value = 82763.5
$ws.Range("L109").Value = 82763.5
$ws.Range("N109").Value = -85537.5
$ws.Range("H122").Value = 3595.9355
$ws.Range("I122").Value = 3403.72
$ws.Range("K122").Value = 10211.16
$ws.Range("M122").Value = -7761.16
$ws.Range("H136").Value = 2691.2808
$ws.Range("I136").Value = 2148.889
$ws.Range("J136").Value = 4725.25
$ws.Range("K136").Value = 6446.667
$ws.Range("L136").Value = 14175.75
$ws.Range("M136").Value = -3896.667
$ws.Range("N136").Value = -19275.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2966.5
$ws.Range("J5").Value = 3888.6667
$ws.Range("L5").Value = 3888.6667
$ws.Range("N5").Value = -4114.6667
$ws.Range("H86").Value = 4645.6665
$ws.Range("I86").Value = 1940
$ws.Range("J86").Value = 5998.5
$ws.Range("K86").Value = 1940
$ws.Range("L86").Value = 5998.5
$ws.Range("M86").Value = -817
$ws.Range("N86").Value = -8244.5
$ws.Range("H89").Value = 4645.6665
$ws.Range("I89").Value = 1940
$ws.Range("J89").Value = 5998.5
$ws.Range("K89").Value = 9700
$ws.Range("L89").Value = 29992.5
$ws.Range("M89").Value = -4084
$ws.Range("N89").Value = -41224.5
$ws.Range("H134").Value = 2904.6316
$ws.Range("I134").Value = 2680.5625
$ws.Range("J134").Value = 4099.6665
$ws.Range("K134").Value = 8041.6875
$ws.Range("L134").Value = 12298.9995
$ws.Range("M134").Value = -5506.6875
$ws.Range("N134").Value = -17368.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3263.158
$ws.Range("I31").Value = 2729.2646
$ws.Range("J31").Value = 4052.3914
$ws.Range("K31").Value = 2729.2646
$ws.Range("L31").Value = 4052.3914
$ws.Range("M31").Value = -2434.2646
$ws.Range("N31").Value = -4642.3914
$ws.Range("H34").Value = 3263.158
$ws.Range("I34").Value = 2729.2646
$ws.Range("J34").Value = 4052.3914
$ws.Range("K34").Value = 2729.2646
$ws.Range("L34").Value = 4052.3914
$ws.Range("M34").Value = -2527.2646
$ws.Range("N34").Value = -4456.3914
$ws.Range("H41").Value = 13570.143
$ws.Range("J41").Value = 24663.666
$ws.Range("L41").Value = 24663.666
$ws.Range("N41").Value = -25519.666
$ws.Range("H58").Value = 2933.9788
$ws.Range("I58").Value = 2858.1162
$ws.Range("K58").Value = 2858.1162
$ws.Range("M58").Value = -2655.1162
$ws.Range("H62").Value = 7715.5
$ws.Range("I62").Value = 3206.75
$ws.Range("J62").Value = 25750.5
$ws.Range("K62").Value = 3206.75
$ws.Range("L62").Value = 25750.5
$ws.Range("M62").Value = -2582.75
$ws.Range("N62").Value = -26998.5
$ws.Range("H65").Value = 7715.5
$ws.Range("I65").Value = 3206.75
$ws.Range("J65").Value = 25750.5
$ws.Range("K65").Value = 16033.75
$ws.Range("L65").Value = 128752.5
$ws.Range("M65").Value = -12913.75
$ws.Range("N65").Value = -134992.5
$ws.Range("H136").Value = 2933.9788
$ws.Range("I136").Value = 2858.1162
$ws.Range("K136").Value = 8574.3486
$ws.Range("M136").Value = -6024.348599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 39446756
$ws.Range("I4").Value = 26890792
$ws.Range("J4").Value = 77899400
$ws.Range("K4").Value = 80672376
$ws.Range("L4").Value = 233698200
$ws.Range("M4").Value = -80672264
$ws.Range("N4").Value = -233698424
$ws.Range("H11").Value = 2001425.5
$ws.Range("I11").Value = 2364381.8
$ws.Range("K11").Value = 7093145.399999999
$ws.Range("M11").Value = -7093005.399999999
$ws.Range("H13").Value = 285.9091
$ws.Range("I13").Value = 279
$ws.Range("J13").Value = 289.85715
$ws.Range("K13").Value = 837
$ws.Range("L13").Value = 869.5714499999999
$ws.Range("M13").Value = -669
$ws.Range("N13").Value = -1205.57145
$ws.Range("H15").Value = 147.66667
$ws.Range("J15").Value = 258.125
$ws.Range("L15").Value = 774.375
$ws.Range("N15").Value = -1054.375
$ws.Range("H17").Value = 388.63635
$ws.Range("J17").Value = 733.3333
$ws.Range("L17").Value = 2199.9999
$ws.Range("N17").Value = -2537.9999
$ws.Range("H26").Value = 3015.2856
$ws.Range("I26").Value = 3015.2856
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 9045.856800000001
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -8757.856800000001
$ws.Range("N26").ClearContents()
$ws.Range("H34").Value = 4506.706
$ws.Range("J34").Value = 5555.2593
$ws.Range("L34").Value = 16665.7779
$ws.Range("N34").Value = -16833.7779
$ws.Range("H37").Value = 39980
$ws.Range("J37").Value = 39980
$ws.Range("L37").Value = 119940
$ws.Range("N37").Value = -120164
$ws.Range("H41").Value = 553.9
$ws.Range("I41").Value = 877.8
$ws.Range("J41").Value = 230
$ws.Range("K41").Value = 2633.4
$ws.Range("L41").Value = 690
$ws.Range("M41").Value = -2295.4
$ws.Range("N41").Value = -1366
$ws.Range("H44").Value = 862.36365
$ws.Range("J44").Value = 618
$ws.Range("L44").Value = 1854
$ws.Range("N44").Value = -2650
$ws.Range("H46").Value = 7143.636
$ws.Range("J46").Value = 4299.4
$ws.Range("L46").Value = 12898.2
$ws.Range("N46").Value = -13080.2
$ws.Range("H56").Value = 7320
$ws.Range("I56").Value = 7320
$ws.Range("K56").Value = 7320
$ws.Range("M56").Value = -6790
$ws.Range("H60").Value = 4534.5386
$ws.Range("I60").Value = 2966
$ws.Range("J60").Value = 4739.1304
$ws.Range("K60").Value = 8898
$ws.Range("L60").Value = 14217.3912
$ws.Range("M60").Value = -8647
$ws.Range("N60").Value = -14719.3912
$ws.Range("H113").Value = 2916.4707
$ws.Range("J113").Value = 2540.0667
$ws.Range("L113").Value = 7620.2001
$ws.Range("N113").Value = -11960.2001
$ws.Range("H122").Value = 2209.9614
$ws.Range("I122").Value = 2241.6667
$ws.Range("J122").Value = 2182.7856
$ws.Range("K122").Value = 20175.0003
$ws.Range("L122").Value = 19645.0704
$ws.Range("M122").Value = -17725.0003
$ws.Range("N122").Value = -24545.0704
$ws.Range("H131").Value = 144142.3
$ws.Range("I131").Value = 2125750
$ws.Range("J131").Value = 2598.8928
$ws.Range("K131").Value = 6377250
$ws.Range("L131").Value = 7796.678400000001
$ws.Range("M131").Value = -6372210
$ws.Range("N131").Value = -17876.6784
$ws.Range("H139").Value = 3594.389
$ws.Range("I139").Value = 2977.8
$ws.Range("K139").Value = 8933.400000000001
$ws.Range("M139").Value = -3793.400000000001
$ws.Range("H140").Value = 2655.9644
$ws.Range("I140").Value = 2655.9644
$ws.Range("K140").Value = 7967.8932
$ws.Range("M140").Value = -2787.8932
$ws.Range("H141").Value = 5928.913
$ws.Range("I141").Value = 5701.7
$ws.Range("K141").Value = 17105.1
$ws.Range("M141").Value = -11925.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 154.54546
$ws.Range("I2").Value = 148.15
$ws.Range("K2").Value = 148.15
$ws.Range("M2").Value = -35.15000000000001
$ws.Range("H24").Value = 27076.924
$ws.Range("J24").Value = 27076.924
$ws.Range("L24").Value = 27076.924
$ws.Range("N24").Value = -27422.924
$ws.Range("H43").Value = 13186.857
$ws.Range("J43").Value = 26009
$ws.Range("L43").Value = 26009
$ws.Range("N43").Value = -26311
$ws.Range("H46").Value = 46255.125
$ws.Range("J46").Value = 42857.145
$ws.Range("L46").Value = 42857.145
$ws.Range("N46").Value = -43169.145
$ws.Range("H97").Value = 900
$ws.Range("I97").Value = 900
$ws.Range("K97").Value = 900
$ws.Range("M97").Value = -404
$ws.Range("H102").Value = 19436.465
$ws.Range("I102").Value = 19971.148
$ws.Range("K102").Value = 19971.148
$ws.Range("M102").Value = -18349.148
$ws.Range("H122").Value = 2293.1892
$ws.Range("I122").Value = 1623.45
$ws.Range("K122").Value = 4870.35
$ws.Range("M122").Value = -2420.35
$ws.Range("H126").Value = 2608.077
$ws.Range("I126").Value = 2410.6365
$ws.Range("K126").Value = 7231.9095
$ws.Range("M126").Value = -4761.9095
$ws.Range("H132").Value = 4459.579
$ws.Range("I132").Value = 3762.2144
$ws.Range("K132").Value = 11286.6432
$ws.Range("M132").Value = -8756.643199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8410
$ws.Range("I7").Value = 4225
$ws.Range("K7").Value = 4225
$ws.Range("M7").Value = -4113
$ws.Range("H16").Value = 4953.75
$ws.Range("I16").Value = 3930.75
$ws.Range("J16").Value = 6999.75
$ws.Range("K16").Value = 3930.75
$ws.Range("L16").Value = 6999.75
$ws.Range("M16").Value = -3760.75
$ws.Range("N16").Value = -7339.75
$ws.Range("H22").Value = 1605.16
$ws.Range("I22").Value = 874.2
$ws.Range("J22").Value = 2092.4666
$ws.Range("K22").Value = 874.2
$ws.Range("L22").Value = 2092.4666
$ws.Range("M22").Value = -579.2
$ws.Range("N22").Value = -2682.4666
$ws.Range("H27").Value = 1605.16
$ws.Range("I27").Value = 874.2
$ws.Range("J27").Value = 2092.4666
$ws.Range("K27").Value = 874.2
$ws.Range("L27").Value = 2092.4666
$ws.Range("M27").Value = -767.2
$ws.Range("N27").Value = -2306.4666
$ws.Range("H38").Value = 54496.7
$ws.Range("J38").Value = 48282.43
$ws.Range("L38").Value = 48282.43
$ws.Range("N38").Value = -49102.43
$ws.Range("H46").Value = 5561.625
$ws.Range("J46").Value = 5561.625
$ws.Range("L46").Value = 5561.625
$ws.Range("N46").Value = -5937.625
$ws.Range("H68").Value = 1773.75
$ws.Range("I68").Value = 1773.75
$ws.Range("K68").Value = 1773.75
$ws.Range("M68").Value = -1024.75
$ws.Range("H71").Value = 1773.75
$ws.Range("I71").Value = 1773.75
$ws.Range("K71").Value = 8868.75
$ws.Range("M71").Value = -5124.75
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H126").Value = 8410
$ws.Range("I126").Value = 4225
$ws.Range("K126").Value = 12675
$ws.Range("M126").Value = -10205
$ws.Range("H132").Value = 2760.52
$ws.Range("I132").Value = 2695.4736
$ws.Range("J132").Value = 2966.5
$ws.Range("K132").Value = 8086.4208
$ws.Range("L132").Value = 8899.5
$ws.Range("M132").Value = -5556.4208
$ws.Range("N132").Value = -13959.5
$ws.Range("H136").Value = 9089.904
$ws.Range("I136").Value = 9431.75
$ws.Range("K136").Value = 28295.25
$ws.Range("M136").Value = -25745.25
$ws.Range("H140").Value = 158333.33
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 158333.33
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 158333.33
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -168693.33

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10900083
$ws.Range("I122").Value = 16711731
$ws.Range("K122").Value = 50135193
$ws.Range("M122").Value = -50132743
$ws.Range("H132").Value = 8068487.5
$ws.Range("I132").Value = 16672278
$ws.Range("J132").Value = 2433.8125
$ws.Range("K132").Value = 50016834
$ws.Range("L132").Value = 7301.4375
$ws.Range("M132").Value = -50014304
$ws.Range("N132").Value = -12361.4375
$ws.Range("H136").Value = 12927.288
$ws.Range("I136").Value = 14432.059
$ws.Range("K136").Value = 43296.177
$ws.Range("M136").Value = -40746.177
